$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# printf(“%d\n”, ageArray[0]);  ->  display(ageArray[0]);
$oldPrintf = [string]::Format("printf({0}%d\n{1}, ageArray[0]);", [char]0x201C, [char]0x201D)
$found1 = $d.Content.Find.Execute($oldPrintf, $true, $false, $false, $false, $false, $true, 1, $false, "display(ageArray[0]);", 2)
Write-Output "Change 1 (printf -> display) applied: $found1"

# --- Change 2 -------------------------------------------------------------
# "... a data type e.g. 4 for int" -> "... a data type or variable e.g. 4 for int"
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("e.g. 4 for int", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $insPoint2 = $d.Range($rng2.Start, $rng2.Start)
    $insPoint2.InsertBefore("or variable ")
}
Write-Output "Change 2 (insert 'or variable ') applied: $found2"

# --- Change 3 ---------------------------------------------------------------
# "...decays into a pointer to the first element" ->
# "...decays into a pointer to the first element. So using sizeof(arr) in the method would only return the size of the pointer instead."
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("decays into a pointer to the first element", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Collapse(0)
    $rng3.InsertAfter(". So using sizeof(arr) in the method would only return the size of the pointer instead.")
}
Write-Output "Change 3 (append sizeof explanation) applied: $found3"
